# Appends the 06-10-2020 COVID19 daily snapshot block (header + 35 states/UTs)
# to the bottom of the COVID19_TIMESERIESDATA sheet, mirroring the existing
# repeating per-day block structure already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRow = 685

# --- Header row (row 685) -------------------------------------------------
# Copy formatting from the most recent existing header block (row 649, the
# last "States/UT" header row already in the sheet) so the new header picks
# up the same bold / centered / thin-bordered look used by every other
# per-day header row, then overwrite with the (identical) header text.
$ws.Range("A649:H649").Copy()
$ws.Range("A685:H685").PasteSpecial(-4122)

$ws.Cells.Item($headerRow, 1).Value = "States/UT"
$ws.Cells.Item($headerRow, 2).Value = "Active Cases"
$ws.Cells.Item($headerRow, 3).Value = "Active Cases Since Yesterday"
$ws.Cells.Item($headerRow, 4).Value = "Recovered Cases"
$ws.Cells.Item($headerRow, 5).Value = "Recovered Cases Since Yesterday"
$ws.Cells.Item($headerRow, 6).Value = "Deceased Cases"
$ws.Cells.Item($headerRow, 7).Value = "Deceased Cases Since Yesterday"
$ws.Cells.Item($headerRow, 8).Value = "Date"

# --- Data rows (686-720), one per state/UT, all dated 06-10-2020 ---------
$data = @(
    @('Andaman and Nicobar Islands', 186, 4, 3659, 10, 54, 1, '06-10-2020'),
    @('Andhra Pradesh', 51060, -3340, 666433, 7558, 6019, 38, '06-10-2020'),
    @('Arunachal Pradesh', 2989, 36, 7775, 198, 19, 1, '06-10-2020'),
    @('Assam', 33467, 143, 153491, 1364, 760, 11, '06-10-2020'),
    @('Bihar', 11523, -272, 176995, 1537, 924, 9, '06-10-2020'),
    @('Chandigarh', 1604, -69, 10797, 199, 177, 3, '06-10-2020'),
    @('Chhattisgarh', 27857, -691, 97067, 3336, 1081, 36, '06-10-2020'),
    @('Dadra and Nagar Haveli and Daman and Diu', 99, -6, 2991, 11, 2, 0, '06-10-2020'),
    @('Delhi', 23080, -1673, 263938, 3588, 5542, 32, '06-10-2020'),
    @('Goa', 4803, -36, 30456, 423, 460, 4, '06-10-2020'),
    @('Gujarat', 16718, -91, 123638, 1405, 3509, 13, '06-10-2020'),
    @('Haryana', 11822, -245, 121596, 1255, 1491, 21, '06-10-2020'),
    @('Himachal Pradesh', 3156, -117, 12653, 292, 224, 7, '06-10-2020'),
    @('Jammu and Kashmir', 14696, -764, 63790, 1386, 1252, 10, '06-10-2020'),
    @('Jharkhand', 10436, -500, 76843, 1312, 747, 4, '06-10-2020'),
    @('Karnataka', 115496, -97, 522846, 7064, 9370, 84, '06-10-2020'),
    @('Kerala', 84958, 379, 149111, 4640, 859, 23, '06-10-2020'),
    @('Ladakh', 1166, 60, 3414, 60, 61, 0, '06-10-2020'),
    @('Madhya Pradesh', 18757, -615, 115878, 2046, 2463, 29, '06-10-2020'),
    @('Maharashtra', 252721, -3001, 1162585, 12982, 38347, 263, '06-10-2020'),
    @('Manipur', 2696, 120, 9334, 129, 75, 1, '06-10-2020'),
    @('Meghalaya', 2217, 8, 4491, 98, 59, 5, '06-10-2020'),
    @('Mizoram', 291, -22, 1837, 30, 0, 0, '06-10-2020'),
    @('Nagaland', 1155, -71, 5422, 113, 17, 0, '06-10-2020'),
    @('Odisha', 28006, -1498, 206400, 4098, 924, 17, '06-10-2020'),
    @('Puducherry', 4513, -274, 24221, 458, 543, 4, '06-10-2020'),
    @('Punjab', 12895, -682, 102648, 1671, 3641, 38, '06-10-2020'),
    @('Rajasthan', 21215, 61, 123421, 2090, 1559, 14, '06-10-2020'),
    @('Sikkim', 598, -51, 2547, 67, 46, 1, '06-10-2020'),
    @('Tamil Nadu', 45881, -239, 569664, 5572, 9846, 62, '06-10-2020'),
    @('Telengana', 26644, -408, 174769, 2381, 1181, 10, '06-10-2020'),
    @('Tripura', 4876, 18, 22131, 255, 301, 2, '06-10-2020'),
    @('Uttarakhand', 8701, -388, 42621, 881, 669, 17, '06-10-2020'),
    @('Uttar Pradesh', 45024, -1361, 366321, 4269, 6092, 63, '06-10-2020'),
    @('West Bengal', 27717, 278, 240707, 3009, 5255, 61, '06-10-2020')
)

$firstDataRow = $headerRow + 1
$lastDataRow = $firstDataRow + $data.Length - 1
$dateRange = $ws.Range("H" + $firstDataRow + ":H" + $lastDataRow)

# Force column H of the new data rows to Text format first so date-shaped
# strings like "06-10-2020" are stored verbatim instead of being parsed
# into a date serial number; the data rows in the source sheet carry no
# explicit style, so clear the formatting back off again once the literal
# text value has been committed.
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $firstDataRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
}

$dateRange.ClearFormats()
